$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 160.64285
$ws.Range("I6").Value = 160.64285
$ws.Range("K6").Value = 481.92855
$ws.Range("M6").Value = -369.92855
$ws.Range("H74").Value = 41672908
$ws.Range("I74").Value = 50004690
$ws.Range("J74").Value = 14000
$ws.Range("K74").Value = 50004690
$ws.Range("L74").Value = 14000
$ws.Range("M74").Value = -50003754
$ws.Range("N74").Value = -15872
$ws.Range("H77").Value = 41672908
$ws.Range("I77").Value = 50004690
$ws.Range("J77").Value = 14000
$ws.Range("K77").Value = 250023450
$ws.Range("L77").Value = 70000
$ws.Range("M77").Value = -250018770
$ws.Range("N77").Value = -79360
$ws.Range("H98").Value = 71435336
$ws.Range("I98").Value = 76929520
$ws.Range("K98").Value = 76929520
$ws.Range("M98").Value = -76928022
$ws.Range("H116").Value = 11909875
$ws.Range("I116").Value = 19233542
$ws.Range("J116").Value = 8915.5
$ws.Range("K116").Value = 19233542
$ws.Range("L116").Value = 8915.5
$ws.Range("M116").Value = -19230100
$ws.Range("N116").Value = -15799.5
$ws.Range("H122").Value = 71435336
$ws.Range("I122").Value = 76929520
$ws.Range("K122").Value = 230788560
$ws.Range("M122").Value = -230786110
$ws.Range("H126").Value = 103963.336
$ws.Range("J126").Value = 103963.336
$ws.Range("L126").Value = 103963.336
$ws.Range("N126").Value = -113843.336
$ws.Range("H132").Value = 1890.375
$ws.Range("I132").Value = 1598.0667
$ws.Range("J132").Value = 2767.3
$ws.Range("K132").Value = 4794.2001
$ws.Range("L132").Value = 8301.900000000001
$ws.Range("M132").Value = -2264.2001
$ws.Range("N132").Value = -13361.9
$ws.Range("H135").Value = 435539.88
$ws.Range("I135").Value = 526949.8
$ws.Range("K135").Value = 4742548.2
$ws.Range("M135").Value = -4740013.2
$ws.Range("H137").Value = 2468.0454
$ws.Range("I137").Value = 2416.2856
$ws.Range("K137").Value = 7248.8568
$ws.Range("M137").Value = -4698.8568
$ws.Range("H138").Value = 6989.567
$ws.Range("J138").Value = 8761.857
$ws.Range("L138").Value = 26285.571
$ws.Range("N138").Value = -36565.571
$ws.Range("H141").Value = 1643
$ws.Range("I141").Value = 1569.8572
$ws.Range("K141").Value = 4709.571599999999
$ws.Range("M141").Value = 470.4284000000007

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5404.886
$ws.Range("I61").Value = 2493.8333
$ws.Range("J61").Value = 11642.857
$ws.Range("K61").Value = 2493.8333
$ws.Range("L61").Value = 11642.857
$ws.Range("M61").Value = -2281.8333
$ws.Range("N61").Value = -12066.857
$ws.Range("H74").Value = 36885.895
$ws.Range("I74").Value = 44551.105
$ws.Range("K74").Value = 44551.105
$ws.Range("M74").Value = -43677.105
$ws.Range("H77").Value = 36885.895
$ws.Range("I77").Value = 44551.105
$ws.Range("K77").Value = 222755.525
$ws.Range("M77").Value = -218387.525
$ws.Range("H136").Value = 5404.886
$ws.Range("I136").Value = 2493.8333
$ws.Range("J136").Value = 11642.857
$ws.Range("K136").Value = 7481.499899999999
$ws.Range("L136").Value = 34928.571
$ws.Range("M136").Value = -4931.499899999999
$ws.Range("N136").Value = -40028.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 68997
$ws.Range("J81").Value = 68997
$ws.Range("L81").Value = 68997
$ws.Range("N81").Value = -71119
$ws.Range("H84").Value = 68997
$ws.Range("J84").Value = 68997
$ws.Range("L84").Value = 206991
$ws.Range("N84").Value = -217599
$ws.Range("H94").Value = 1381.1072
$ws.Range("I94").Value = 929.17645
$ws.Range("J94").Value = 2079.5454
$ws.Range("K94").Value = 929.17645
$ws.Range("L94").Value = 2079.5454
$ws.Range("M94").Value = -478.17645
$ws.Range("N94").Value = -2981.5454
$ws.Range("H134").Value = 4802.482
$ws.Range("I134").Value = 1856.8182
$ws.Range("K134").Value = 5570.4546
$ws.Range("M134").Value = -3035.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4683.1113
$ws.Range("I16").Value = 2928
$ws.Range("J16").Value = 6312.857
$ws.Range("K16").Value = 2928
$ws.Range("L16").Value = 6312.857
$ws.Range("M16").Value = -2641
$ws.Range("N16").Value = -6886.857
$ws.Range("H31").Value = 8636.780000000001
$ws.Range("J31").Value = 14736.05
$ws.Range("L31").Value = 14736.05
$ws.Range("N31").Value = -15326.05
$ws.Range("H34").Value = 8636.780000000001
$ws.Range("J34").Value = 14736.05
$ws.Range("L34").Value = 14736.05
$ws.Range("N34").Value = -15140.05
$ws.Range("H58").Value = 11117624
$ws.Range("I58").Value = 26318690
$ws.Range("K58").Value = 26318690
$ws.Range("M58").Value = -26318487
$ws.Range("H105").Value = 3108668.8
$ws.Range("I105").Value = 3969799
$ws.Range("K105").Value = 3969799
$ws.Range("M105").Value = -3968052
$ws.Range("H113").Value = 4683.1113
$ws.Range("I113").Value = 2928
$ws.Range("J113").Value = 6312.857
$ws.Range("K113").Value = 2928
$ws.Range("L113").Value = 6312.857
$ws.Range("M113").Value = -758
$ws.Range("N113").Value = -10652.857
$ws.Range("H132").Value = 4709.364
$ws.Range("I132").Value = 1931.7826
$ws.Range("J132").Value = 7751.476
$ws.Range("K132").Value = 5795.3478
$ws.Range("L132").Value = 23254.428
$ws.Range("M132").Value = -3265.3478
$ws.Range("N132").Value = -28314.428
$ws.Range("H134").Value = 5672.2666
$ws.Range("I134").Value = 3401.7097
$ws.Range("K134").Value = 10205.1291
$ws.Range("M134").Value = -7670.1291
$ws.Range("H136").Value = 11117624
$ws.Range("I136").Value = 26318690
$ws.Range("K136").Value = 78956070
$ws.Range("M136").Value = -78953520

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 366.66666
$ws.Range("H12").Value = 4167774.2
$ws.Range("I12").Value = 1147.75
$ws.Range("J12").Value = 6251087.5
$ws.Range("K12").Value = 3443.25
$ws.Range("L12").Value = 18753262.5
$ws.Range("M12").Value = -3270.25
$ws.Range("N12").Value = -18753608.5
$ws.Range("H14").Value = 13889522
$ws.Range("I14").Value = 13889522
$ws.Range("K14").Value = 41668566
$ws.Range("M14").Value = -41668393
$ws.Range("H34").Value = 4025.8076
$ws.Range("I34").Value = 182.5
$ws.Range("J34").Value = 5178.8
$ws.Range("K34").Value = 547.5
$ws.Range("L34").Value = 15536.4
$ws.Range("M34").Value = -463.5
$ws.Range("N34").Value = -15704.4
$ws.Range("H76").Value = 500001500
$ws.Range("I76").Value = 500001500
$ws.Range("K76").Value = 1500004500
$ws.Range("M76").Value = -1500004117
$ws.Range("H79").Value = 500001500
$ws.Range("I79").Value = 500001500
$ws.Range("K79").Value = 1500004500
$ws.Range("M79").Value = -1500003174
$ws.Range("H80").Value = 250001000
$ws.Range("I80").Value = 250001000
$ws.Range("J80").Value = 250001000
$ws.Range("K80").Value = 750003000
$ws.Range("L80").Value = 750003000
$ws.Range("M80").Value = -750002064
$ws.Range("N80").Value = -750004872
$ws.Range("H83").Value = 250001000
$ws.Range("I83").Value = 250001000
$ws.Range("J83").Value = 250001000
$ws.Range("K83").Value = 2250009000
$ws.Range("L83").Value = 2250009000
$ws.Range("M83").Value = -2250004320
$ws.Range("N83").Value = -2250018360
$ws.Range("H92").Value = 25642358
$ws.Range("J92").Value = 25642358
$ws.Range("L92").Value = 76927074
$ws.Range("N92").Value = -76929570
$ws.Range("H107").Value = 33335600
$ws.Range("I107").Value = 700
$ws.Range("J107").Value = 50003050
$ws.Range("K107").Value = 2100
$ws.Range("L107").Value = 150009150
$ws.Range("M107").Value = -180
$ws.Range("N107").Value = -150012990

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1601320
$ws.Range("J107").Value = 2000
$ws.Range("L107").Value = 2000
$ws.Range("N107").Value = -5840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 26501.25
$ws.Range("J23").Value = 34999
$ws.Range("L23").Value = 34999
$ws.Range("N23").Value = -35459
$ws.Range("H55").Value = 66667336
$ws.Range("I55").Value = 1000000000
$ws.Range("K55").Value = 1000000000
$ws.Range("M55").Value = -999999827
$ws.Range("H93").Value = 5053.647
$ws.Range("I93").Value = 6722.4443
$ws.Range("K93").Value = 6722.4443
$ws.Range("M93").Value = -5474.4443
$ws.Range("H132").Value = 7942577
$ws.Range("I132").Value = 13161785
$ws.Range("J132").Value = 9380.92
$ws.Range("K132").Value = 39485355
$ws.Range("L132").Value = 28142.76
$ws.Range("M132").Value = -39482825
$ws.Range("N132").Value = -33202.76

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10648490
$ws.Range("I132").Value = 14710541
$ws.Range("K132").Value = 44131623
$ws.Range("M132").Value = -44129093
